$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'298.06"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'2.10%"
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'42.19"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'4.87%"
$ws.Range("E3").ClearFormats()
$ws.Range("E4").Value = "'0.01%"
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'0.07520"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'3.16%"
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'1.591"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'2.24%"
$ws.Range("E6").ClearFormats()
$ws.Range("D7").Value = "'0.9257"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'-0.14%"
$ws.Range("E7").ClearFormats()
$ws.Range("E8").Value = "'1.71%"
$ws.Range("E8").ClearFormats()
$ws.Range("E9").Value = "'2.14%"
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'0.1833"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'4.01%"
$ws.Range("E10").ClearFormats()
$ws.Range("D11").Value = "'0.08938"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'2.02%"
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'0.04132"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'-5.25%"
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'0.1048"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'-0.37%"
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'0.001278"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'0.30%"
$ws.Range("E14").ClearFormats()
$ws.Range("D15").Value = "'0.005871"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'-2.40%"
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'3.338"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'-0.12%"
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'4.367"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'1.99%"
$ws.Range("E17").ClearFormats()
$ws.Range("E18").Value = "'1.36%"
$ws.Range("E18").ClearFormats()
$ws.Range("D19").Value = "'8.342"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'5.18%"
$ws.Range("E19").ClearFormats()
$ws.Range("E20").Value = "'-2.78%"
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'0.3104"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'11.88%"
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'0.04083"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'3.92%"
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'0.001266"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'0.27%"
$ws.Range("E23").ClearFormats()
$ws.Range("E24").Value = "'5.92%"
$ws.Range("E24").ClearFormats()
$ws.Range("D25").Value = "'0.0001301"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'8.27%"
$ws.Range("E25").ClearFormats()
$ws.Range("D38").Value = "'0.02398"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'3.93%"
$ws.Range("E38").ClearFormats()
$ws.Range("D39").Value = "'0.05224"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "'2.94%"
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = "'0.006791"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'25.41%"
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = "'0.007765"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'-1.13%"
$ws.Range("E41").ClearFormats()
$ws.Range("D42").Value = "'0.1324"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'3.06%"
$ws.Range("E42").ClearFormats()
$ws.Range("D43").Value = "'0.007417"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'0.26%"
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'0.007113"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "'-2.29%"
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'0.2991"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'3.08%"
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = "'0.00006575"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'6.28%"
$ws.Range("E46").ClearFormats()
$ws.Range("E47").Value = "'-0.06%"
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'0.05487"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'13.00%"
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'0.004204"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'0.05%"
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'-0.06%"
$ws.Range("E50").ClearFormats()
$ws.Range("E51").Value = "'-0.06%"
$ws.Range("E51").ClearFormats()
